# Update the dSF (column F) values to reflect the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    5  = 2
    6  = 3
    11 = 0
    14 = 5
    15 = 2
    20 = -1
    21 = 3
    22 = 6
    23 = -5
    24 = 1
    25 = 3
    27 = -1
    31 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
